$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / rich-text strings: new report week ---
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Crime statistics table updates (rows 14-29) ---
# Step 1: write new values (numbers as numbers; "0"/"***.*" placeholders as text
#         via a leading apostrophe so Excel does not coerce them to numbers).
$ws.Range("C14").Value = "'0"
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 66.666666666666
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 11.111111111111
$ws.Range("N15").Value = -64.285714285714
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = -65.217391304347
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = -48.275862068965
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = -47.058823529411
$ws.Range("N16").Value = -83.516483516483
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 62.5
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 134
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = -14.649681528662
$ws.Range("L17").Value = -2.898550724637
$ws.Range("M17").Value = 65.432098765432
$ws.Range("N17").Value = -49.433962264150
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -69.565217391304
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = -33.75
$ws.Range("L18").Value = 1.923076923076
$ws.Range("M18").Value = 60.606060606060
$ws.Range("N18").Value = -82.214765100671
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 3.333333333333
$ws.Range("I19").Value = 102
$ws.Range("J19").Value = 105
$ws.Range("K19").Value = -2.857142857142
$ws.Range("L19").Value = 12.087912087912
$ws.Range("M19").Value = 17.241379310344
$ws.Range("N19").Value = -15
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 37.5
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -6.896551724137
$ws.Range("L20").Value = 92.857142857142
$ws.Range("M20").Value = 170
$ws.Range("N20").Value = -70.329670329670
$ws.Range("C21").Value = 30
$ws.Range("E21").Value = -3.225806451612
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -29.710144927536
$ws.Range("I21").Value = 374
$ws.Range("J21").Value = 470
$ws.Range("K21").Value = -20.425531914893
$ws.Range("L21").Value = 10.324483775811
$ws.Range("M21").Value = 22.222222222222
$ws.Range("N21").Value = -65.561694290976
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 12.5
$ws.Range("I23").Value = 67
$ws.Range("J23").Value = 65
$ws.Range("K23").Value = 3.076923076923
$ws.Range("L23").Value = 6.349206349206
$ws.Range("M23").Value = 39.583333333333
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -26.315789473684
$ws.Range("F24").Value = 63
$ws.Range("H24").Value = -4.545454545454
$ws.Range("I24").Value = 256
$ws.Range("J24").Value = 265
$ws.Range("K24").Value = -3.396226415094
$ws.Range("L24").Value = 13.274336283185
$ws.Range("M24").Value = 27.363184079602
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -38.461538461538
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = -37.878787878787
$ws.Range("I25").Value = 168
$ws.Range("J25").Value = 182
$ws.Range("K25").Value = -7.692307692307
$ws.Range("L25").Value = 18.309859154929
$ws.Range("M25").Value = -32.258064516129
$ws.Range("C26").Value = "'0"
$ws.Range("F26").Value = 6
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 15
$ws.Range("K26").Value = 36.363636363636
$ws.Range("L26").Value = 36.363636363636
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 37.5
$ws.Range("I27").Value = 19
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -17.391304347826
$ws.Range("L27").Value = -26.923076923076
$ws.Range("C28").Value = "'0"
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = -11.111111111111
$ws.Range("N28").Value = -75
$ws.Range("C29").Value = "'0"
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = -12.5
$ws.Range("L29").Value = -30
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -75

# Step 2: a few cells flip between a numeric count and the text placeholder "0"
# or "***.*". Re-apply the number/text cell style from a same-column neighbour
# that already carries the correct style, so the style id matches the new cell type.
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
